# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values are recalculated/regenerated and written
# back into the worksheet for rows 2-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 5
    4  = 7
    5  = 4
    6  = 2
    7  = 5
    8  = 5
    9  = 3
    10 = 4
    11 = 6
    12 = 5
    13 = 5
    14 = 4
    15 = 10
    16 = 6
    17 = 1
    18 = 6
    19 = 3
    20 = 8
    21 = 5
    22 = 5
    23 = 5
    24 = 5
    25 = 7
    26 = 2
    27 = 3
    28 = 1
    29 = 4
    30 = 6
    31 = 10
    32 = 5
    33 = 4
    34 = 7
    35 = 7
    36 = 3
    37 = 3
    38 = 3
    39 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
